$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-5 from 2023-09-16 (45185)
# to 2023-10-05 (45204), keeping existing cell formatting/style intact.
$ws.Range("C2:C5").Value = 45204
